{"js": "// Update the first diary entry (\"EXAMPLE STYLE 1\"): change its date from\n// 10.9.2018 to 07.06.2021 and replace the body paragraph's narrative text\n// with the new text about setting up git / choosing VS Code.\n\nconst body = context.document.body;\n\n// Locate the target paragraph uniquely via its old narrative text (there are\n// two \"10.9.2018\" dates in the document - this search pins us to the right\n// one, scoped to its own paragraph, and leaves the other date untouched).\nconst anchorHits = body.search(\"checked the general information\", { matchCase: true });\nanchorHits.load(\"items\");\nawait context.sync();\n\nif (anchorHits.items.length === 0) {\n  throw new Error(\"Could not find the diary paragraph to update.\");\n}\n\nconst targetParagraph = anchorHits.items[0].paragraphs.getFirst();\ntargetParagraph.load(\"text\");\nawait context.sync();\n\n// 1) Swap the date run's text in place (keeps the run's own formatting).\nconst dateHits = targetParagraph.search(\"10.9.2018\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\n\nif (dateHits.items.length > 0) {\n  dateHits.items[0].insertText(\"07.06.2021\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Replace the old narrative (\"I checked the general information ... intro\n//    to GIT. \") with the new narrative. The trailing single space run after\n//    this chunk (kept as-is) is intentionally left out of the match so it\n//    survives untouched, exactly like in the source diff.\nconst oldNarrative =\n  \"I checked the general information and understood the main focus of the \" +\n  \"course, which is to find my passion as a software developer and create \" +\n  \"a unique project to represent my skills. I chose frontend module \" +\n  \"because it was the most interesting project offered. I\\u2019ve also \" +\n  \"tried to set up my environment, but I could not decide which code \" +\n  \"editor I would like to use.  I learned to set up a git repository and \" +\n  \"did my first commit, everything went smoothly after I clicked the \" +\n  \"banner to watch intro to GIT. \";\n\nconst newNarrative =\n  \"I set up a git repository, read the instructions for completing the \" +\n  \"course. I chose VS Code as my editor since I\\u2019ve used it before \" +\n  \"and have it set up with all the extensions.  \";\n\nconst narrativeHits = targetParagraph.search(oldNarrative, { matchCase: true });\nnarrativeHits.load(\"items\");\nawait context.sync();\n\nif (narrativeHits.items.length > 0) {\n  narrativeHits.items[0].insertText(newNarrative, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the first diary entry (\"EXAMPLE STYLE 1\"): change its date from\n# 10.9.2018 to 07.06.2021 and replace the body paragraph's narrative text\n# with the new text about setting up git / choosing VS Code.\n\n$d = $word.ActiveDocument\n\n# Locate the target paragraph uniquely via its old narrative text (there are\n# two \"10.9.2018\" dates in the document - this anchors us to the right one\n# and leaves the other date untouched).\n$anchor = $d.Content\n$anchor.Find.ClearFormatting()\n$anchor.Find.Text = \"checked the general information\"\n$anchor.Find.MatchCase = $true\n$anchor.Find.Wrap = 0\n$found = $anchor.Find.Execute()\nif (-not $found) {\n  throw \"Could not find the diary paragraph to update.\"\n}\n\n$targetPara = $anchor.Paragraphs(1).Range\n\n# 1) Swap the date text in place (scoped to the target paragraph only).\n$dateRange = $targetPara.Duplicate\n$dateRange.Find.ClearFormatting()\n$dateRange.Find.Text = \"10.9.2018\"\n$dateRange.Find.MatchCase = $true\n$dateRange.Find.Wrap = 0\nif ($dateRange.Find.Execute()) {\n  $dateRange.Text = \"07.06.2021\"\n}\n\n# 2) Replace the old narrative (\"I checked the general information ... intro\n#    to GIT. \") with the new narrative. The trailing single space run after\n#    this chunk (kept as-is) is intentionally left out of the match so it\n#    survives untouched, exactly like in the source diff.\n$oldNarrative = \"I checked the general information and understood the main focus of the course, which is to find my passion as a software developer and create a unique project to represent my skills. I chose frontend module because it was the most interesting project offered. I\" + [char]0x2019 + \"ve also tried to set up my environment, but I could not decide which code editor I would like to use.  I learned to set up a git repository and did my first commit, everything went smoothly after I clicked the banner to watch intro to GIT. \"\n$newNarrative = \"I set up a git repository, read the instructions for completing the course. I chose VS Code as my editor since I\" + [char]0x2019 + \"ve used it before and have it set up with all the extensions.  \"\n\n$narrativeRange = $targetPara.Duplicate\n$narrativeRange.Find.ClearFormatting()\n$narrativeRange.Find.Text = $oldNarrative\n$narrativeRange.Find.MatchCase = $true\n$narrativeRange.Find.Wrap = 0\nif ($narrativeRange.Find.Execute()) {\n  $narrativeRange.Text = $newNarrative\n}\n"}
